# Scheduled-runner data refresh: re-pull currentAveragePrice(NQ/HQ) market
# data from Universalis and recompute the dependent Leve profit columns
# (H:N) for the affected crafting-leve rows across the job sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 227534.1
$ws.Range("I103").Value = 357288.56
$ws.Range("J103").Value = 463.75
$ws.Range("K103").Value = 1071865.68
$ws.Range("L103").Value = 1391.25
$ws.Range("M103").Value = -1071279.68
$ws.Range("N103").Value = -2563.25
$ws.Range("H111").Value = 2095.5173
$ws.Range("I111").Value = 1244.0869
$ws.Range("J111").Value = 5359.3335
$ws.Range("K111").Value = 3732.2607
$ws.Range("L111").Value = 16078.0005
$ws.Range("M111").Value = -665.2606999999998
$ws.Range("N111").Value = -22212.0005
$ws.Range("H125").Value = 900
$ws.Range("I125").Value = 100
$ws.Range("J125").Value = 1300
$ws.Range("K125").Value = 900
$ws.Range("L125").Value = 11700
$ws.Range("M125").Value = 1560
$ws.Range("N125").Value = -16620
$ws.Range("H137").Value = 31185.314
$ws.Range("I137").Value = 2920.1052
$ws.Range("J137").Value = 64750.25
$ws.Range("K137").Value = 8760.3156
$ws.Range("L137").Value = 194250.75
$ws.Range("M137").Value = -6210.3156
$ws.Range("N137").Value = -199350.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18631.555
$ws.Range("I32").Value = 22918.125
$ws.Range("K32").Value = 22918.125
$ws.Range("M32").Value = -22631.125
$ws.Range("H61").Value = 324093.12
$ws.Range("I61").Value = 475763.2
$ws.Range("K61").Value = 475763.2
$ws.Range("M61").Value = -475551.2
$ws.Range("H74").Value = 2487.8928
$ws.Range("I74").Value = 3322.5625
$ws.Range("J74").Value = 1375
$ws.Range("K74").Value = 3322.5625
$ws.Range("L74").Value = 1375
$ws.Range("M74").Value = -2448.5625
$ws.Range("N74").Value = -3123
$ws.Range("H77").Value = 2487.8928
$ws.Range("I77").Value = 3322.5625
$ws.Range("J77").Value = 1375
$ws.Range("K77").Value = 16612.8125
$ws.Range("L77").Value = 6875
$ws.Range("M77").Value = -12244.8125
$ws.Range("N77").Value = -15611
$ws.Range("H122").Value = 2329.9412
$ws.Range("I122").Value = 2252
$ws.Range("J122").Value = 2517
$ws.Range("K122").Value = 6756
$ws.Range("L122").Value = 7551
$ws.Range("M122").Value = -4306
$ws.Range("N122").Value = -12451
$ws.Range("H132").Value = 10812.164
$ws.Range("I132").Value = 1779.9111
$ws.Range("J132").Value = 51457.3
$ws.Range("K132").Value = 5339.7333
$ws.Range("L132").Value = 154371.9
$ws.Range("M132").Value = -2809.7333
$ws.Range("N132").Value = -159431.9
$ws.Range("H136").Value = 324093.12
$ws.Range("I136").Value = 475763.2
$ws.Range("K136").Value = 1427289.6
$ws.Range("M136").Value = -1424739.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 34894.16
$ws.Range("I134").Value = 43052.56
$ws.Range("J134").Value = 900.8333
$ws.Range("K134").Value = 129157.68
$ws.Range("L134").Value = 2702.4999
$ws.Range("M134").Value = -126622.68
$ws.Range("N134").Value = -7772.4999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10000.868
$ws.Range("I31").Value = 15532.827
$ws.Range("J31").Value = 3316.4167
$ws.Range("K31").Value = 15532.827
$ws.Range("L31").Value = 3316.4167
$ws.Range("M31").Value = -15237.827
$ws.Range("N31").Value = -3906.4167
$ws.Range("H34").Value = 10000.868
$ws.Range("I34").Value = 15532.827
$ws.Range("J34").Value = 3316.4167
$ws.Range("K34").Value = 15532.827
$ws.Range("L34").Value = 3316.4167
$ws.Range("M34").Value = -15330.827
$ws.Range("N34").Value = -3720.4167
$ws.Range("H134").Value = 4962.2915
$ws.Range("I134").Value = 691.8421
$ws.Range("J134").Value = 21190
$ws.Range("K134").Value = 2075.5263
$ws.Range("L134").Value = 63570
$ws.Range("M134").Value = 459.4737
$ws.Range("N134").Value = -68640

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4026.7715
$ws.Range("J68").Value = 4439.9355
$ws.Range("L68").Value = 13319.8065
$ws.Range("N68").Value = -14941.8065
$ws.Range("H71").Value = 4026.7715
$ws.Range("J71").Value = 4439.9355
$ws.Range("L71").Value = 39959.4195
$ws.Range("N71").Value = -48071.4195
$ws.Range("H120").Value = 14298
$ws.Range("I120").Value = 6696.6665
$ws.Range("J120").Value = 19999
$ws.Range("K120").Value = 20089.9995
$ws.Range("L120").Value = 59997
$ws.Range("M120").Value = -15251.9995
$ws.Range("N120").Value = -69673
$ws.Range("H131").Value = 132408.17
$ws.Range("J131").Value = 147897.52
$ws.Range("L131").Value = 443692.5599999999
$ws.Range("N131").Value = -453772.5599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 28625
$ws.Range("I24").Value = 100000
$ws.Range("J24").Value = 4833.3335
$ws.Range("K24").Value = 100000
$ws.Range("L24").Value = 4833.3335
$ws.Range("M24").Value = -99827
$ws.Range("N24").Value = -5179.3335
$ws.Range("H102").Value = 7602.2
$ws.Range("I102").Value = 8877.75
$ws.Range("K102").Value = 8877.75
$ws.Range("M102").Value = -7255.75
$ws.Range("H113").Value = 3619.8
$ws.Range("I113").Value = 2039.8
$ws.Range("K113").Value = 2039.8
$ws.Range("M113").Value = 130.2
$ws.Range("H122").Value = 11000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 33000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -37900
$ws.Range("H126").Value = 6354.48
$ws.Range("I126").Value = 5647.625
$ws.Range("J126").Value = 7611.1113
$ws.Range("K126").Value = 16942.875
$ws.Range("L126").Value = 22833.3339
$ws.Range("M126").Value = -14472.875
$ws.Range("N126").Value = -27773.3339
$ws.Range("H132").Value = 39984.3
$ws.Range("I132").Value = 55586.05
$ws.Range("K132").Value = 166758.15
$ws.Range("M132").Value = -164228.15

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5193
$ws.Range("I7").Value = 5350.2354
$ws.Range("J7").Value = 4524.75
$ws.Range("K7").Value = 5350.2354
$ws.Range("L7").Value = 4524.75
$ws.Range("M7").Value = -5238.2354
$ws.Range("N7").Value = -4748.75
$ws.Range("H40").Value = 98089.586
$ws.Range("I40").Value = 189552.5
$ws.Range("J40").Value = 6626.6665
$ws.Range("K40").Value = 189552.5
$ws.Range("L40").Value = 6626.6665
$ws.Range("M40").Value = -189416.5
$ws.Range("N40").Value = -6898.6665
$ws.Range("H46").Value = 1045.0714
$ws.Range("I46").Value = 794.25
$ws.Range("J46").Value = 2550
$ws.Range("K46").Value = 794.25
$ws.Range("L46").Value = 2550
$ws.Range("M46").Value = -606.25
$ws.Range("N46").Value = -2926
$ws.Range("H55").Value = 338.25
$ws.Range("I55").Value = 173.33333
$ws.Range("J55").Value = 437.2
$ws.Range("K55").Value = 173.33333
$ws.Range("L55").Value = 437.2
$ws.Range("M55").Value = -0.3333299999999895
$ws.Range("N55").Value = -783.2
$ws.Range("H61").Value = 4053
$ws.Range("I61").Value = 2011.5
$ws.Range("J61").Value = 6775
$ws.Range("K61").Value = 2011.5
$ws.Range("L61").Value = 6775
$ws.Range("M61").Value = -1809.5
$ws.Range("N61").Value = -7179
$ws.Range("H113").Value = 4053
$ws.Range("I113").Value = 2011.5
$ws.Range("J113").Value = 6775
$ws.Range("K113").Value = 2011.5
$ws.Range("L113").Value = 6775
$ws.Range("M113").Value = 158.5
$ws.Range("N113").Value = -11115
$ws.Range("H122").Value = 3117.389
$ws.Range("I122").Value = 2701.1428
$ws.Range("K122").Value = 8103.428400000001
$ws.Range("M122").Value = -5653.428400000001
$ws.Range("H126").Value = 5193
$ws.Range("I126").Value = 5350.2354
$ws.Range("J126").Value = 4524.75
$ws.Range("K126").Value = 16050.7062
$ws.Range("L126").Value = 13574.25
$ws.Range("M126").Value = -13580.7062
$ws.Range("N126").Value = -18514.25
$ws.Range("H132").Value = 1476.58
$ws.Range("I132").Value = 1106.2368
$ws.Range("J132").Value = 2649.3333
$ws.Range("K132").Value = 3318.7104
$ws.Range("L132").Value = 7947.999899999999
$ws.Range("M132").Value = -788.7103999999999
$ws.Range("N132").Value = -13007.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1643.75
$ws.Range("I107").Value = 468.875
$ws.Range("J107").Value = 2818.625
$ws.Range("K107").Value = 1406.625
$ws.Range("L107").Value = 8455.875
$ws.Range("M107").Value = 513.375
$ws.Range("N107").Value = -12295.875
$ws.Range("H113").Value = 2703914.8
$ws.Range("I113").Value = 1863.3334
$ws.Range("J113").Value = 6756992
$ws.Range("K113").Value = 5590.0002
$ws.Range("L113").Value = 20270976
$ws.Range("M113").Value = -3420.0002
$ws.Range("N113").Value = -20275316
$ws.Range("H122").Value = 1999.8667
$ws.Range("I122").Value = 1837.8334
$ws.Range("K122").Value = 5513.5002
$ws.Range("M122").Value = -3063.5002
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 9000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 1532.238
$ws.Range("I132").Value = 1360.7878
$ws.Range("J132").Value = 2160.889
$ws.Range("K132").Value = 4082.3634
$ws.Range("L132").Value = 6482.667
$ws.Range("M132").Value = -1552.3634
$ws.Range("N132").Value = -11542.667
$ws.Range("H136").Value = 1493.2084
$ws.Range("I136").Value = 852.75
$ws.Range("K136").Value = 2558.25
$ws.Range("M136").Value = -8.25
